$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '27.519.67'
$ws.Range('E2').Value2 = '  -0.18%  '
$ws.Range('D3').Value2 = '1.617.93'
$ws.Range('E3').Value2 = '  -1.40%  '
$ws.Range('E4').Value2 = '  +0.08%  '
$ws.Range('D6').Value2 = '''0.525'
$ws.Range('E6').Value2 = '  -2.08%  '
$ws.Range('E7').Value2 = '  +0.07%  '
$ws.Range('D8').Value2 = '''22.68'
$ws.Range('E8').Value2 = '  -0.78%  '
$ws.Range('E9').Value2 = '  +2.38%  '
$ws.Range('E10').Value2 = '  +0.21%  '
$ws.Range('D11').Value2 = '''0.0886'
$ws.Range('E11').Value2 = '  -0.31%  '
$ws.Range('D12').Value2 = '1.847.74'
$ws.Range('E12').Value2 = '  -1.32%  '
$ws.Range('D13').Value2 = '1.619.18'
$ws.Range('E13').Value2 = '  -1.15%  '
$ws.Range('E14').Value2 = '  -0.22%  '
$ws.Range('E15').Value2 = '  -1.65%  '
$ws.Range('E16').Value2 = '  +1.51%  '
$ws.Range('D17').Value2 = '27.528.31'
$ws.Range('E17').Value2 = '  -0.01%  '
$ws.Range('D18').Value2 = '''229.82'
$ws.Range('E18').Value2 = '  +0.74%  '
$ws.Range('D20').Value2 = '''7.52'
$ws.Range('E20').Value2 = '  -1.19%  '
$ws.Range('E21').Value2 = '  -0.02%  '
$ws.Range('E22').Value2 = '  -0.03%  '
$ws.Range('D23').Value2 = '''10.11'
$ws.Range('E23').Value2 = '  +1.35%  '
$ws.Range('E24').Value2 = '  +6.98%  '
$ws.Range('D25').Value2 = '''149.59'
$ws.Range('E25').Value2 = '  +0.29%  '
$ws.Range('E26').Value2 = '  -1.24%  '
$ws.Range('E27').Value2 = '  +0.01%  '
$ws.Range('D28').Value2 = '''6.79'
$ws.Range('E28').Value2 = '  -2.21%  '
$ws.Range('D29').Value2 = '''15.53'
$ws.Range('E29').Value2 = '  -0.18%  '
$ws.Range('E30').Value2 = '  -0.28%  '
$ws.Range('D31').Value2 = '''0.0481'
$ws.Range('E31').Value2 = '  -0.69%  '
$ws.Range('E32').Value2 = '  -0.88%  '
$ws.Range('D33').Value2 = '1.441.08'
$ws.Range('E33').Value2 = '  +1.01%  '
$ws.Range('E34').Value2 = '  -3.33%  '
$ws.Range('E35').Value2 = '  -3.11%  '
$ws.Range('D36').Value2 = '''2.33'
$ws.Range('E36').Value2 = '  -0.34%  '
$ws.Range('D37').Value2 = '''0.933'
$ws.Range('E37').Value2 = '  +3.48%  '
$ws.Range('E38').Value2 = '  -2.12%  '
$ws.Range('E39').Value2 = '  +0.21%  '
$ws.Range('D40').Value2 = '''0.859'
$ws.Range('E40').Value2 = '  -1.72%  '
$ws.Range('D41').Value2 = '''69.10'
$ws.Range('E41').Value2 = '  +6.30%  '
$ws.Range('E42').Value2 = '  -0.01%  '
$ws.Range('E43').Value2 = '  -3.09%  '
$ws.Range('D44').Value2 = '''2.45'
$ws.Range('E44').Value2 = '  -0.48%  '
$ws.Range('E45').Value2 = '  -1.76%  '
$ws.Range('E46').Value2 = '  -2.00%  '
$ws.Range('D47').Value2 = '1.758.04'
$ws.Range('E47').Value2 = '  -1.34%  '
$ws.Range('E48').Value2 = '  +0.41%  '
$ws.Range('D49').Value2 = '''86.27'
$ws.Range('E49').Value2 = '  +0.07%  '
$ws.Range('E50').Value2 = '  -1.82%  '
$ws.Range('D51').Value2 = '''0.0995'
$ws.Range('E51').Value2 = '  +1.33%  '
